$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the value "NA" in G4 and G5 (matches existing shared string "NA" used elsewhere, e.g. B2:B5/C2:C5)
$ws.Range("G4").Value = "NA"
$ws.Range("G5").Value = "NA"

# Update the active selection to G5 (was C14)
$ws.Range("G5").Select()
